# Add "Data (dd/mm/yyyy)" column (G) to the activation log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header ---
$ws.Range("G5").Value = "Data (dd/mm/yyyy)"

# --- Data: same dates as column E (Date of Activation), re-displayed dd/mm/yyyy ---
$ws.Range("G6").Value = 44563
$ws.Range("G7").Value = 44635
$ws.Range("G8").Value = 44655
$ws.Range("G9").Value = 44701
$ws.Range("G10").Value = 44706
$ws.Range("G11").Value = 44718
$ws.Range("G12").Value = 44749
$ws.Range("G13").Value = 44781
$ws.Range("G14").Value = 44813
$ws.Range("G15").Value = 44844
$ws.Range("G16").Value = 44876
$ws.Range("G17").Value = 44907
$ws.Range("G18").Value = 44908
$ws.Range("G19").Value = 44909
$ws.Range("G20").Value = 44941
$ws.Range("G21").Value = 44973
$ws.Range("G22").Value = 44974
$ws.Range("G23").Value = 45003
$ws.Range("G24").Value = 45004
$ws.Range("G25").Value = 45036
$ws.Range("G26").Value = 45037
$ws.Range("G27").Value = 45038
$ws.Range("G28").Value = 45069
$ws.Range("G29").Value = 45070
$ws.Range("G30").Value = 45071
$ws.Range("G31").Value = 45103
$ws.Range("G32").Value = 45104
$ws.Range("G33").Value = 45135
$ws.Range("G34").Value = 45136
$ws.Range("G35").Value = 45168

# --- Formatting: header + data rows only (don't touch untouched rows 1-4) ---
$ws.Range("G5:G35").Font.Name = "Book Antiqua"
$ws.Range("G5:G35").Font.Size = 12
$ws.Range("G5:G35").NumberFormat = "dd/mm/yyyy"
$ws.Columns("G").ColumnWidth = 22.1

# --- View: scrolled down a few more rows, new selection ---
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L24").Select()
